# TC03_Canine_Filter_PrimDisSite-Lung.xlsx - "10 icdc scripts for jenkins"
#
# The FilesTab Neo4j query (cell B4 on the "startup" sheet) is updated to
# drop the `File Type` and `Breed` columns from its RETURN clause, and the
# sheet's active selection moves from D2 to B4 to reflect the cell that was
# edited.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$newQuery = @'

MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
 MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
WHERE diag.primary_disease_site IN ['Lung']
WITH DISTINCT f, parent, c, demo, diag, s
RETURN coalesce(f.file_name, '') AS `File Name`, 
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

$ws.Range("B4").Value = $newQuery

# Reflect the new active cell/selection state on the sheet.
[void]$ws.Range("B4").Select()
